$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the D2 header text (shared string) ---
$ws.Range("D2").Value2 = "f'(x) = f(x+dx)-f(x) / f(x)"

# --- Wrap the header text and widen column D / heighten row 2 to fit ---
$ws.Range("D2").WrapText = $true
$ws.Columns.Item(4).ColumnWidth = 24.5
$ws.Rows.Item(2).RowHeight = 36

# --- Add a new data row (row 12) by copying formatting down from row 11 ---
$ws.Range("B11:D11").Copy()
$ws.Range("B12:D12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("B12").Value2 = 10

# --- C column: re-enter as one shared formula across C3:C12 ---
$ws.Range("C3:C12").Formula = "=B3^2"

# --- D column: replace the "double" formula with a forward-difference formula ---
$ws.Range("D3").ClearContents()
$ws.Range("D4:D11").FormulaR1C1 = "=R[1]C[-1]-RC[-1]"

# --- Selection / active cell ---
$ws.Range("D3").Select()
